# =====================================================================
# Update with Correct Forecast output
# Renames Sheet1 -> "Sales vs PO", inserts an "Order Week" column, and
# adds three new sheets: "Weekly Growth", "Volume Insights", "Prediction Info"
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# --- Insert a new "Order Week" column at C (shifts PO_Requested_Qty C -> D) ---
$ws1.Columns.Item(3).Insert()
$ws1.Cells.Item(1,3).Value = "Order Week"

# --- Rewrite the data rows: col A (new ds), col C (Order Week = old ds), col D (PO_Requested_Qty) ---
$ws1.Cells.Item(2,1).Value = 45494
$ws1.Cells.Item(2,3).Value = 45488
$ws1.Cells.Item(2,4).Value = 0
$ws1.Cells.Item(3,1).Value = 45508
$ws1.Cells.Item(3,3).Value = 45502
$ws1.Cells.Item(3,4).Value = 0
$ws1.Cells.Item(4,1).Value = 45515
$ws1.Cells.Item(4,3).Value = 45509
$ws1.Cells.Item(4,4).Value = 0
$ws1.Cells.Item(5,1).Value = 45522
$ws1.Cells.Item(5,3).Value = 45516
$ws1.Cells.Item(5,4).Value = 0
$ws1.Cells.Item(6,1).Value = 45529
$ws1.Cells.Item(6,3).Value = 45523
$ws1.Cells.Item(6,4).Value = 0
$ws1.Cells.Item(7,1).Value = 45536
$ws1.Cells.Item(7,3).Value = 45530
$ws1.Cells.Item(7,4).Value = 0
$ws1.Cells.Item(8,1).Value = 45543
$ws1.Cells.Item(8,3).Value = 45537
$ws1.Cells.Item(8,4).Value = 0
$ws1.Cells.Item(9,1).Value = 45550
$ws1.Cells.Item(9,3).Value = 45544
$ws1.Cells.Item(9,4).Value = 0
$ws1.Cells.Item(10,1).Value = 45557
$ws1.Cells.Item(10,3).Value = 45551
$ws1.Cells.Item(10,4).Value = 0
$ws1.Cells.Item(11,1).Value = 45564
$ws1.Cells.Item(11,3).Value = 45558
$ws1.Cells.Item(11,4).Value = 0
$ws1.Cells.Item(12,1).Value = 45571
$ws1.Cells.Item(12,3).Value = 45565
$ws1.Cells.Item(12,4).Value = 0
$ws1.Cells.Item(13,1).Value = 45578
$ws1.Cells.Item(13,3).Value = 45572
$ws1.Cells.Item(13,4).Value = 0
$ws1.Cells.Item(14,1).Value = 45585
$ws1.Cells.Item(14,3).Value = 45579
$ws1.Cells.Item(14,4).Value = 0
$ws1.Cells.Item(15,1).Value = 45592
$ws1.Cells.Item(15,3).Value = 45586
$ws1.Cells.Item(15,4).Value = 0
$ws1.Cells.Item(16,1).Value = 45599
$ws1.Cells.Item(16,3).Value = 45593
$ws1.Cells.Item(16,4).Value = 0
$ws1.Cells.Item(17,1).Value = 45606
$ws1.Cells.Item(17,3).Value = 45600
$ws1.Cells.Item(17,4).Value = 0
$ws1.Cells.Item(18,1).Value = 45613
$ws1.Cells.Item(18,3).Value = 45607
$ws1.Cells.Item(18,4).Value = 0
$ws1.Cells.Item(19,1).Value = 45620
$ws1.Cells.Item(19,3).Value = 45614
$ws1.Cells.Item(19,4).Value = 0
$ws1.Cells.Item(20,1).Value = 45627
$ws1.Cells.Item(20,3).Value = 45621
$ws1.Cells.Item(20,4).Value = 0
$ws1.Cells.Item(21,1).Value = 45634
$ws1.Cells.Item(21,3).Value = 45628
$ws1.Cells.Item(21,4).Value = 0
$ws1.Cells.Item(22,1).Value = 45641
$ws1.Cells.Item(22,3).Value = 45635
$ws1.Cells.Item(22,4).Value = 0
$ws1.Cells.Item(23,1).Value = 45648
$ws1.Cells.Item(23,3).Value = 45642
$ws1.Cells.Item(23,4).Value = 0
$ws1.Cells.Item(24,1).Value = 45655
$ws1.Cells.Item(24,3).Value = 45649
$ws1.Cells.Item(24,4).Value = 0

# --- Apply the date number format (copied from column A) onto the new Order Week column ---
$ws1.Range("A2:A24").Copy()
$ws1.Range("C2:C24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add Sheet 2: "Weekly Growth" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"
$ws2.Cells.Item(1,1).Value = "ds"
$ws2.Cells.Item(1,2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1,3).Value = "Growth%"
$ws2.Cells.Item(2,1).Value = 45495
$ws2.Cells.Item(2,2).Value = 20
$ws2.Cells.Item(2,3).Value = 0
$ws2.Cells.Item(3,1).Value = 45509
$ws2.Cells.Item(3,2).Value = 20
$ws2.Cells.Item(3,3).Value = 0
$ws2.Cells.Item(4,1).Value = 45523
$ws2.Cells.Item(4,2).Value = 20
$ws2.Cells.Item(4,3).Value = 0
$ws2.Cells.Item(5,1).Value = 45530
$ws2.Cells.Item(5,2).Value = 20
$ws2.Cells.Item(5,3).Value = 0
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("A2:A5").Copy()
$ws2.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add Sheet 3: "Volume Insights" ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"
$ws3.Cells.Item(1,1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1,2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1,3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1,4).Value = "Min_PO_Quantity"
$ws3.Cells.Item(2,1).Value = 80
$ws3.Cells.Item(2,2).Value = 20
$ws3.Cells.Item(2,3).Value = 20
$ws3.Cells.Item(2,4).Value = 20
$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add Sheet 4: "Prediction Info" ---
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"
$ws4.Cells.Item(1,1).Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Cells.Item(2,1).Value = 20
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Restore Sheet1 as the active sheet ---
$ws1.Activate()
